# Applies the cryptos price/volume refresh + FirstDigitalUSD/PolygonEcosystemToken row swap
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.343.27"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "'2.661.95"
$ws.Range("E3").Value = "  +3.37%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'611.04"
$ws.Range("E5").Value = "  +3.76%  "

$ws.Range("D6").Value = "'143.44"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").Value = "'2.660.34"
$ws.Range("E9").Value = "  +3.34%  "

$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("D11").Value = "'5.61"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").Value = "'0.362"
$ws.Range("E13").Value = "  +3.11%  "

$ws.Range("D14").Value = "'27.27"
$ws.Range("E14").Value = "  +0.09%  "

$ws.Range("D15").Value = "'3.139.10"
$ws.Range("E15").Value = "  +3.47%  "

$ws.Range("D16").Value = "'63.188.90"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").Value = "'0.0000144"
$ws.Range("E17").Value = "  -1.00%  "

$ws.Range("D18").Value = "'2.668.38"
$ws.Range("E18").Value = "  +3.83%  "

$ws.Range("D19").Value = "'11.45"
$ws.Range("E19").Value = "  +3.52%  "

$ws.Range("D20").Value = "'341.53"
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("E21").Value = "  +1.87%  "

$ws.Range("D22").Value = "'6.87"
$ws.Range("E22").Value = "  +3.52%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").Value = "'66.99"
$ws.Range("E24").Value = "  -1.30%  "

$ws.Range("E25").Value = "  +1.54%  "

$ws.Range("D26").Value = "'1.53"
$ws.Range("E26").Value = "  -2.53%  "

$ws.Range("D27").Value = "'8.64"
$ws.Range("E27").Value = "  +4.88%  "

$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("D29").Value = "'545.48"
$ws.Range("E29").Value = "  +15.43%  "

$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").Value = "'7.84"
$ws.Range("E31").Value = "  -1.19%  "

$ws.Range("D32").Value = "'2.05"
$ws.Range("E32").Value = "  +5.22%  "

$ws.Range("E33").Value = "  +6.59%  "

$ws.Range("D34").Value = "'0.0₃0805"
$ws.Range("E34").Value = "  +0.36%  "

$ws.Range("D35").Value = "'173.38"
$ws.Range("E35").Value = "  -1.65%  "

$ws.Range("D36").Value = "'5.19"
$ws.Range("E36").Value = "  +13.99%  "

$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.406"
$ws.Range("E37").Value = "  +1.38%  "

$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").Value = "'19.15"
$ws.Range("E39").Value = "  +1.64%  "

$ws.Range("D40").Value = "'1.85"
$ws.Range("E40").Value = "  +8.92%  "

$ws.Range("D41").Value = "'176.42"
$ws.Range("E41").Value = "  +11.87%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("E43").Value = "  +1.64%  "

$ws.Range("D44").Value = "'22.21"
$ws.Range("E44").Value = "  +4.07%  "

$ws.Range("E45").Value = "  +6.73%  "

$ws.Range("D46").Value = "'0.633"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("E47").Value = "  +1.35%  "

$ws.Range("D48").Value = "'0.0962"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D49").Value = "'18.72"
$ws.Range("E49").Value = "  +3.26%  "

$ws.Range("D50").Value = "'1.74"
$ws.Range("E50").Value = "  +3.92%  "
